$d = $word.ActiveDocument

# 1. Split the "Critério" run into two sentences separated by a manual line break.
$d.Content.Find.Execute(
    "A Nota final (NF) será calculada da seguinte maneira: NF = 50%(P1) + 50%(P2)Cada docente responsável usará seu próprio critério na aplicação de trabalhos.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "A Nota final (NF) será calculada da seguinte maneira: NF = 50%(P1) + 50%(P2)^lCada docente responsável usará seu próprio critério na aplicação de trabalhos.",
    2)

# 2. Split the Bibliografia run into one entry per line (manual line breaks),
#    with a blank line before "Bibliografia complementar:".
$d.Content.Find.Execute(
    "SMITH, J.M.; VAN NESS, H.C.; ABBOTT, M.M.; SWIHART, M.T. Introduction to Chemical Engineering Thermodynamics. 9th ed. Editora McGraw Hill, 2022.SANDLER, S.I., Chemical, Biochemical, and Engineering Thermodynamics, 5th ed., Editora John Wiley & Sons, 2020 TERRON, L. R. Termodinâmica Química Aplicada. 1 ed. Editora Manole Ltda, 2009.Bibliografia complementar:MATSOUKAS, T. Fundamentos de Termodinâmica para Engenharia Química. 1 ed. LTC Editora, 2016.TAVARES, F.W.; SEGTOVICH, I.S.V.; MEDEIROS, F.A. Termodinâmica na Engenharia Química. 1ra ed. LTC Editora, 2023.BALZISHER, R.E.; SAMUELS M.R.; ELIASSEN J.D. Termodinámica Química para Ingenieros. Prentice-Hall Inc., 1974.KORETSKY, M. D. Termodinâmica para Engenharia Química. 1 ed. LTC Editora, 2007.MORAN, M. I.; SHAPIRO, H. N.; BOETTNER, D.D.; BAILEY, M.B.  Fundamentals of Engineering Thermodynamics. 9th. Editora John Wiley & Sons, 2018. BORGNAKKE, C.; SONNTAG, R.E. Fundamentos da Termodinâmica Clássica. 8th ed. Editora Blucher, 2013",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "SMITH, J.M.; VAN NESS, H.C.; ABBOTT, M.M.; SWIHART, M.T. Introduction to Chemical Engineering Thermodynamics. 9th ed. Editora McGraw Hill, 2022.^lSANDLER, S.I., Chemical, Biochemical, and Engineering Thermodynamics, 5th ed., Editora John Wiley & Sons, 2020 ^lTERRON, L. R. Termodinâmica Química Aplicada. 1 ed. Editora Manole Ltda, 2009.^l^lBibliografia complementar:^lMATSOUKAS, T. Fundamentos de Termodinâmica para Engenharia Química. 1 ed. LTC Editora, 2016.^lTAVARES, F.W.; SEGTOVICH, I.S.V.; MEDEIROS, F.A. Termodinâmica na Engenharia Química. 1ra ed. LTC Editora, 2023.^lBALZISHER, R.E.; SAMUELS M.R.; ELIASSEN J.D. Termodinámica Química para Ingenieros. Prentice-Hall Inc., 1974.^lKORETSKY, M. D. Termodinâmica para Engenharia Química. 1 ed. LTC Editora, 2007.^lMORAN, M. I.; SHAPIRO, H. N.; BOETTNER, D.D.; BAILEY, M.B.  Fundamentals of Engineering Thermodynamics. 9th. Editora John Wiley & Sons, 2018. ^lBORGNAKKE, C.; SONNTAG, R.E. Fundamentos da Termodinâmica Clássica. 8th ed. Editora Blucher, 2013",
    2)
